$wb = $excel.ActiveWorkbook

$wsCodes = $wb.Worksheets.Item("Codes")
$wsRates = $wb.Worksheets.Item("Rates")

# Add the new "Services" header to the Rates sheet (column G), matching
# the formatting already used by the other header cells in row 1.
$wsRates.Range("F1").Copy() | Out-Null
$wsRates.Range("G1").PasteSpecial(-4122) | Out-Null
$wsRates.Range("G1").Value = "Services"
$excel.CutCopyMode = 0

# Update the selection on the Codes sheet
$wsCodes.Range("E1").Select() | Out-Null

# Update the selection on the Rates sheet and make it the active tab
$wsRates.Activate()
$wsRates.Range("L20").Select() | Out-Null

$wb.Save()
